$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for match rows 136, 137, 138 (sheet rows) is rotated between each other:
#   new row136 = old row138 data (B, E:AD)
#   new row137 = old row136 data (B, E:AD)
#   new row138 = old row137 data (B, E:AD)
# Column A (row index), C (Div), D (Date) stay unchanged.

$data = @{
    136 = @{
        B = 8229445
        E = "NK Vodice"
        F = "NK Zadar"
        G = 1
        H = 1
        I = 0
        J = 1
        K = "D"
        L = 7
        M = 4.5
        N = 1.333
        O = 7
        P = 4.5
        Q = 1.333
        R = 1.5
        S = 1.9
        T = 1.9
        U = 3
        V = 1.825
        W = 1.975
        X = -1
        Y = 3.5
        Z = -1
        AA = 0.8999999999999999
        AB = -1
        AC = -1
        AD = 0.9750000000000001
    }
    137 = @{
        B = 8229444
        E = "NK Neretva"
        F = "Zmaj Makarska"
        G = 1
        H = 0
        I = 1
        J = 0
        K = "H"
        L = 1.727
        M = 3.75
        N = 3.75
        O = 1.5
        P = 4.2
        Q = 4.75
        R = -1
        S = 1.8
        T = 2
        U = 3
        V = 1.875
        W = 1.925
        X = 0.5
        Y = -1
        Z = -1
        AA = 0
        AB = 0
        AC = -1
        AD = 0.925
    }
    138 = @{
        B = 8229446
        E = "RNK Split"
        F = "NK Omis"
        G = 2
        H = 1
        I = 1
        J = 1
        K = "H"
        L = 2.1
        M = 3.4
        N = 2.9
        O = 1.8
        P = 3.6
        Q = 3.6
        R = -0.5
        S = 1.85
        T = 1.95
        U = 2.75
        V = 1.8
        W = 2
        X = 0.8
        Y = -1
        Z = -1
        AA = 0.8500000000000001
        AB = -1
        AC = 0.4
        AD = -0.5
    }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
